# Clean up and commenting
#
# Rows 3..17 (the data rows between the "2175/150/5.5" first record and the
# "1100/55/4" later record) get reshuffled into a new order - same 15
# records, just re-sequenced. Columns F/G (EstateType/DistributionType) are
# constant across these rows so they don't need touching.
#
# Strategy:
#  - Numeric columns (A:D, H:AU) are snapshotted into plain value arrays and
#    written back with Range.Value - this naturally leaves a source blank
#    (e.g. a row with no ConstructionYear) blank at the destination too,
#    with no leftover placeholder cell.
#  - Column E (ZipCode) is text that happens to look numeric ("97074" etc.),
#    so it has to move via Range.Copy (which preserves the shared-string
#    text type) staged through a scratch area, since Copy reads live
#    worksheet data rather than a value we can reorder in memory.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 3
$lastDataRow = 17
$rowCount = $lastDataRow - $firstDataRow + 1
$scratchStartRow = 100 # well below the used range (A1:AU20)

# destOffset[i] (0-based i -> row firstDataRow+i) takes its content from the
# row that was originally at (firstDataRow + order[i] - 1), order being
# 1-based offsets within the 15-row block.
$order = @(6, 5, 11, 7, 3, 2, 4, 9, 1, 12, 10, 8, 14, 15, 13)

# --- Snapshot the numeric columns (A:D and H:AU) for all 15 rows ---
$snapAD = @()
$snapHAU = @()
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstDataRow + $i
    $snapAD += , ($ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 4)).Value())
    $snapHAU += , ($ws.Range($ws.Cells.Item($r, 8), $ws.Cells.Item($r, 47)).Value())
}

# --- Stage column E (ZipCode text) through a scratch area via Copy, which
#     keeps it a shared-string text cell instead of coercing to a number ---
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $firstDataRow + $i
    $scratchRow = $scratchStartRow + $i
    $ws.Cells.Item($r, 5).Copy($ws.Cells.Item($scratchRow, 5))
}

# --- Write back A:D and H:AU in the new order ---
for ($i = 0; $i -lt $rowCount; $i++) {
    $destRow = $firstDataRow + $i
    $srcOffset = $order[$i] - 1
    $ws.Range($ws.Cells.Item($destRow, 1), $ws.Cells.Item($destRow, 4)).Value = $snapAD[$srcOffset]
    $ws.Range($ws.Cells.Item($destRow, 8), $ws.Cells.Item($destRow, 47)).Value = $snapHAU[$srcOffset]
}

# --- Write back column E (ZipCode) in the new order, from the scratch area ---
for ($i = 0; $i -lt $rowCount; $i++) {
    $destRow = $firstDataRow + $i
    $srcOffset = $order[$i] - 1
    $scratchRow = $scratchStartRow + $srcOffset
    $ws.Cells.Item($scratchRow, 5).Copy($ws.Cells.Item($destRow, 5))
}

# --- Wipe the scratch area ---
$ws.Range($ws.Cells.Item($scratchStartRow, 5), $ws.Cells.Item($scratchStartRow + $rowCount - 1, 5)).ClearContents()
